$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Add the new worksheet after the existing one.
$newSheet = $wb.Worksheets.Add($null, $wb.Worksheets.Item($wb.Worksheets.Count))
$newSheet.Name = "Sheet1"

# Copy the daily-data block (header + 31 days) from the original sheet
# (A9:K40) into the new sheet starting at A1 (values + number formats).
$ws1.Range("A9:K40").Copy()
$newSheet.Range("A1").PasteSpecial(-4104)

# A second pass with xlPasteFormats reuses the existing cell styles
# (s="1" header style, s="2" bordered/wrapped data style) instead of
# creating duplicate style entries.
$ws1.Range("A9:K40").Copy()
$newSheet.Range("A1").PasteSpecial(-4122)

$excel.CutCopyMode = 0

# Data rows wrap at the (narrower) default column width on the new
# sheet, which grows each row to two lines.
$newSheet.Rows("2:32").RowHeight = 28.8

[void]$newSheet.Range("A1:K32").Select()
